$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Item starts with C use:" row's Dim Template cell (C4) held a stale
# template filename (F-825-247M ... LLIF Dimension Measure Rev 21_final.xlsx)
# that is no longer a valid/available template. Clear it out (cell becomes
# blank, and the now-unreferenced shared string is dropped on save).
$ws.Range("C4").Value = ""

# Let the row reflow now that it no longer needs to wrap a long filename.
$ws.Rows.Item(4).AutoFit()

# Move the active selection to the cell that was just edited.
$ws.Range("C4").Select()
